# "pedoe pick 3 problems"
# Append three new LeetCode "Tree" problems (picked by Pedoe) to the bottom
# of the "Easy" worksheet table: rows 37, 38 and 39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the three new rows with the formatting used by the existing table rows
# (row 36 is the last data row: bold-styled numeric/string cells in columns
# A, B, D, F, G, matching the "PIC" lookup style used throughout the sheet).
$ws.Range("A36:G36").Copy($ws.Range("A37:G39"))

# Row 37: LeetCode 513 - Find Bottom Left Tree Value (Medium)
$ws.Range("A37").Value = 513
$ws.Range("B37").Value = "Find Bottom Left Tree Value"
$ws.Range("C37").Value = "Tree"
$ws.Range("D37").Value = "Pedoe"
$ws.Range("F37").Value = "Medium"
$ws.Range("G37").Value = "Javascript"

# Row 38: LeetCode 671 - Second Minimum Node in a Binary Tree (Easy)
$ws.Range("A38").Value = 671
$ws.Range("B38").Value = "Second Minimum Node in a Binary Tree"
$ws.Range("C38").Value = "Tree"
$ws.Range("D38").Value = "Pedoe"
$ws.Range("F38").Value = "Easy"
$ws.Range("G38").Value = "Javascript"

# Row 39: LeetCode 530 - Minimum Absolute Difference in BST (Easy)
$ws.Range("A39").Value = 530
$ws.Range("B39").Value = "Minimum Absolute Difference in BST"
$ws.Range("C39").Value = "Tree"
$ws.Range("D39").Value = "Pedoe"
$ws.Range("F39").Value = "Easy"
$ws.Range("G39").Value = "Javascript"

# None of the new rows use column E ("Comment"); drop the blank cells the
# copy above introduced so the rows stay sparse like their neighbours.
$ws.Range("E37:E39").ClearContents()

# Row 37's Difficulty cell needs the "Medium" style (bold/blue), which
# differs from the "Easy" style broadcast by the row-36 copy above; borrow
# the correct formatting from an existing "Medium" cell (F33), then restore
# the value.
$ws.Range("F33").Copy($ws.Range("F37"))
$ws.Range("F37").Value = "Medium"

# Match the author's final selection/cursor position on the sheet.
$ws.Range("C39").Select() | Out-Null
